# [Silverfox] Animset 컬럼 삭제
#
# The "Npc" worksheet had a column "animset" (column D) that is no
# longer needed. Delete the whole column; Excel shifts every column to
# its right (E:J -> D:I) one position to the left to close the gap.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Npc")
$ws.Activate()

# Remove column D ("animset") entirely.
$ws.Columns.Item(4).Delete()

# Leave the selection where the user ended up after the edit.
$ws.Range("G6").Select()
